$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a quiz "Marksheet" that (due to a float-handling bug) was
# duplicating / mis-splitting the answer key across three answer blocks
# (A:B, D:E, G:H) and never actually recorded what the student answered.
# The fix:
#   * Collapses the question list back down to the real 28 questions
#     (25 in block 1 [A:B], 3 in block 2 [D:E]); block 3 [G:H] disappears.
#   * Fills in the "Student Ans" columns (A and D) with what the student
#     actually answered, colouring them green/red/black to mark
#     correct / wrong / not-attempted, exactly like the existing key
#     cells already in the sheet.
#   * Fixes the summary box (rows 10-12) with the recomputed numbers.
# ---------------------------------------------------------------------------

# Reference cells that already carry the three "graded" styles so we can
# clone their formatting (font colour, borders, alignment, ...) without
# creating brand-new style entries.
$refTitle     = $ws.Range("A9")   # mtitleStyle  (bold, centered, bordered)
$refCorrect   = $ws.Range("B10")  # correctStyle (green)
$refIncorrect = $ws.Range("C10")  # incorrectStyle (red)

function Set-Answer($addr, $value, $kind) {
    $cell = $ws.Range($addr)
    $cell.Value = $value
    if ($kind -eq "correct") {
        $refCorrect.Copy()
        $cell.PasteSpecial(-4122)
    } elseif ($kind -eq "wrong") {
        $refIncorrect.Copy()
        $cell.PasteSpecial(-4122)
    }
    # "none" / not-attempted cells already use the default normalStyle
    # formatting, so nothing else needs to change for them.
}

# --- Row labels (A10:A12) become bold/centered like the header row (A9) ---
$refTitle.Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# --- Summary box: Right / Wrong / Not-Attempt / Max ------------------------
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

# --- Marking scheme: marks for correct / wrong (now a real number) --------
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# --- Totals -----------------------------------------------------------------
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "61/112"

# --- Answer block 1 (questions 1-25, columns A = student, B = correct) ----
Set-Answer "A16" "Option A" "correct"
Set-Answer "A17" "Option D" "correct"
Set-Answer "A18" "Option B" "correct"
Set-Answer "A19" "Option C" "correct"
Set-Answer "A20" "Option B" "correct"
Set-Answer "A21" "Option C" "correct"
Set-Answer "A22" "Option D" "correct"
Set-Answer "A25" "Option A" "correct"
Set-Answer "A26" "Option A" "wrong"
Set-Answer "A27" "Option A" "correct"
Set-Answer "A30" "Option B" "correct"
Set-Answer "A31" "Option D" "correct"
Set-Answer "A32" "Option C" "correct"
Set-Answer "A33" "Option D" "correct"
Set-Answer "A34" "Option B" "correct"
Set-Answer "A36" "Option D" "wrong"
Set-Answer "A38" "Option A" "correct"

# --- Answer block 2 (questions 26-28, columns D = student, E = correct) ---
# Only the first 3 rows of this block remain; the rest of the block and all
# of block 3 are removed below.
Set-Answer "D16" "Option A" "correct"
Set-Answer "D18" "Option B" "wrong"

# --- Remove the now-unused part of answer block 2 (rows 19-40) -------------
$ws.Range("D19:E40").Clear()

# --- Remove answer block 3 entirely (columns G:H, rows 15-40) --------------
$ws.Range("G15:H40").Clear()
